# Journal de travail (jnr) weekly update:
#  - fill in this week's entries (the block that starts at row 25, date
#    17.11, previously left blank as a template)
#  - insert one extra detail row into that block (was 5 rows, now 6) to fit
#    an extra task line
#  - keep the named print area in sync with the now-taller sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Insert a row before the current row 29 (the last/bottom row of the
# 25-29 week block). This pushes that bottom row down to 30 and every row
# below it shifts down by one, exactly like "right click row 29 > Insert".
$ws.Rows(29).Insert()

# Re-use the formatting of the row that just got pushed down (row 30, the
# bottom row of the week block) for the freshly inserted row 29, so the new
# row keeps the same borders/alignment/number-formats as the rest of the
# block instead of Excel's bare default style.
$ws.Range("A30:E30").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Rows(29).RowHeight = 18
$excel.CutCopyMode = $false

# --- Fill the week block that starts at row 25 (dates/tasks for this week) ---

# Row 25: header line of the week (date cell already holds =B18+7)
$ws.Range("A25").Value = "Retard"
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = "J'ai eu 10 minutes de retard"

# Row 26
$ws.Range("A26").Value = "Analyse"
$ws.Range("C26").Value = 30
$ws.Range("D26").Value = "J'ai modifié certaine UserStorie"

# Row 27
$ws.Range("A27").Value = "Conception"
$ws.Range("C27").Value = 60
$ws.Range("D27").Value = "J'ai mis en place les déplacement classique avce w,a,s,d,"

# Row 28
$ws.Range("A28").Value = "Conception"
$ws.Range("C28").Value = 60
$ws.Range("D28").Value = "J'ai fait en sorte que le joueur ralentise en s'arrêtant pendant une demi seconde"

# Row 29 (newly inserted row)
$ws.Range("A29").Value = "Conception"
$ws.Range("C29").Value = 10
$ws.Range("D29").Value = "J'ai mis à jour mon git et le jnr "

# Row 30 (was row 29 before the insert)
$ws.Range("A30").Value = "Analyse"
$ws.Range("C30").Value = 10
$ws.Range("D30").Value = "J'ai commencé à regarder comment faire le tire du joueur"

# --- Keep the named print area in sync with the now-taller sheet ---
$ws.PageSetup.PrintArea = '$A$1:$E$68'

$ws.Range("D30").Select()

Write-Output "edit applied"
